$wb = $excel.ActiveWorkbook

# --- Sheet "level": rewrite the per-row track markers (columns A-D) ---
$level = $wb.Worksheets.Item("level")

$level.Range("C1").Value = $null
$level.Range("A1").Value = 1

$level.Range("B2").Value = $null
$level.Range("D2").Value = $null
$level.Range("C2").Value = 1

$level.Range("A3").Value = 1
$level.Range("D3").Value = 1

$level.Range("A4").Value = $null
$level.Range("C4").Value = $null

$level.Range("B5").Value = 1
$level.Range("D5").Value = $null

$level.Range("C7").Value = 1

$level.Range("B8").Value = $null
$level.Range("D8").Value = $null

$level.Range("C9").Value = $null
$level.Range("D9").Value = 7

$level.Range("B12").Value = $null
$level.Range("D12").Value = $null

$level.Range("A15").Value = $null
$level.Range("C15").Value = $null

$level.Range("B17").Value = $null
$level.Range("D17").Value = $null

$level.Range("O5").Select()

# --- Sheet "enemies": append a new enemy entry (row 8) ---
$enemies = $wb.Worksheets.Item("enemies")

$enemies.Range("A8").Value = 7
$enemies.Range("B8").Value = 0
$enemies.Range("C8").Value = 4
$enemies.Range("D8").Value = 0
$enemies.Range("E8").Value = 0.5
$enemies.Range("F8").Value = "EnemyPrefabs/Bullet Enemies/Neo Fly/Neo Fly"

$enemies.Range("A8").HorizontalAlignment = -4108
$enemies.Range("A8").VerticalAlignment = -4108

$enemies.Range("E11").Select()

$level.Select()
